# Update the "Reference ID" value in D2 of the Furniture sheet.
# Original: "41655678"  ->  New: "30677524"
#
# The value is textual (a reference-id code, not a numeric quantity), so a
# helper cell is formatted as Text ("@") before the value is written and the
# result is copied/pasted into D2 - this avoids Excel's automatic
# "numeric-looking string -> number" coercion while keeping D2's own
# style/number-format untouched (matching how the cell was originally
# authored, with its shared-string value rendered through the sheet's
# default style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell far away from the used range as a staging area.
$staging = $ws.Range("Z1")
$staging.NumberFormat = "@"
$staging.Value = "30677524"

$staging.Copy()
$ws.Range("D2").PasteSpecial(-4163)  # xlPasteValues

# Remove the scratch column entirely so it leaves no trace in the sheet.
$ws.Columns.Item(26).Delete()
